# Apply the "adding csv to so" edit: turn the single Soul/currency row
# into a small item table (Soul, Moonstone, Plastic Flower) and trim the
# sheet down to the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Trim the sheet: drop the extra rows (6-15) and extra columns (D:E) ---
$ws.Range("A6:E15").EntireRow.Delete()
$ws.Range("D1:E5").EntireColumn.Delete()

# --- New data rows -----------------------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Moonstone"
$ws.Range("C2").Value = "Material for upgrading weapon"

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Plastic Flower"
$ws.Range("C3").Value = "Delicacy flower, a gift which have not been given"

# --- Row heights (rows 1-3 grew to fit the wrapped description text) ---
$ws.Rows.Item(1).RowHeight = 45.6
$ws.Rows.Item(2).RowHeight = 45.6
$ws.Rows.Item(3).RowHeight = 45.6

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 54.666666666666664

# --- Alignment cleanup ---------------------------------------------------
# Column C description cells: keep vertical centering, drop the explicit
# left horizontal alignment (back to General) that the template had.
$ws.Range("C1").HorizontalAlignment = 1
$ws.Range("C2").HorizontalAlignment = 1
$ws.Range("C3").HorizontalAlignment = 1

# Row 1 keeps the wrapped description style already on C1; rows 2 & 3 use
# a plain vertically centered style (no wrap).
$ws.Range("C2").VerticalAlignment = -4108
$ws.Range("C3").VerticalAlignment = -4108

# --- Selection -------------------------------------------------------
$ws.Range("B3").Select() | Out-Null

Write-Host "edit applied"
